$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.592.24"
$ws.Range("E2").Value = "  +6.50%  "

$ws.Range("D3").Value = "2.039.95"
$ws.Range("E3").Value = "  +3.24%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.87"
$ws.Range("E5").Value = "  +4.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.646"
$ws.Range("E6").Value = "  +2.63%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "66.33"
$ws.Range("E7").Value = "  +18.95%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.48"
$ws.Range("E9").Value = "  +0.34%  "

$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.374"
$ws.Range("E10").Value = "  +5.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0752"
$ws.Range("E11").Value = "  +4.00%  "

$ws.Range("E12").Value = "  +0.97%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.903"
$ws.Range("E13").Value = "  +1.73%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.06"
$ws.Range("E14").Value = "  +6.51%  "

$ws.Range("D15").Value = "2.340.38"
$ws.Range("E15").Value = "  +3.16%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.59"
$ws.Range("E16").Value = "  +7.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.74"
$ws.Range("E17").Value = "  +22.31%  "

$ws.Range("D18").Value = "2.049.71"
$ws.Range("E18").Value = "  +3.46%  "

$ws.Range("D19").Value = "37.535.07"
$ws.Range("E19").Value = "  +6.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.26"
$ws.Range("E20").Value = "  +5.17%  "

$ws.Range("D21").Value = "0.0₃0872"
$ws.Range("E21").Value = "  +5.03%  "

$ws.Range("E22").Value = "  +7.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.85"
$ws.Range("E23").Value = "  +2.53%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.68"
$ws.Range("E24").Value = "  +19.70%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.24%  "

$ws.Range("E26").Value = "  +5.74%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.55"
$ws.Range("E27").Value = "  +6.25%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.10"
$ws.Range("E28").Value = "  +1.45%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.88"
$ws.Range("E29").Value = "  +2.92%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.21"
$ws.Range("E30").Value = "  +10.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.121"
$ws.Range("E31").Value = "  +3.19%  "

$ws.Range("E32").Value = "  +7.27%  "

$ws.Range("E33").Value = "  +24.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.72"
$ws.Range("E34").Value = "  +11.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0611"
$ws.Range("E35").Value = "  +5.18%  "

$ws.Range("E36").Value = "  +8.61%  "

$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.08"
$ws.Range("E38").Value = "  +26.19%  "

$ws.Range("E39").Value = "  +1.89%  "

$ws.Range("E40").Value = "  +17.34%  "

$ws.Range("E41").Value = "  +4.30%  "

$ws.Range("E42").Value = "  +5.06%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.73"
$ws.Range("E43").Value = "  +21.98%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0218"
$ws.Range("E44").Value = "  +5.57%  "

$ws.Range("E45").Value = "  +6.17%  "

$ws.Range("E46").Value = "  +8.97%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.92"
$ws.Range("E47").Value = "  +10.23%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "94.94"
$ws.Range("E48").Value = "  +5.49%  "

$ws.Range("D49").Value = "1.424.76"
$ws.Range("E49").Value = "  +5.73%  "

$ws.Range("E50").Value = "  +2.71%  "

$ws.Range("E51").Value = "  +4.48%  "

